$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = 4
$ws.Range("B24").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C24").Value = "Los Lagos"
$ws.Range("D24").Value = (Get-Date -Year 2022 -Month 11 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D24").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E24").Value = 10
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100107
$ws.Range("H24").Value = "Otros"
$ws.Range("I24").Value = 100107002
$ws.Range("J24").Value = "Chirimoya"
$ws.Range("K24").Value = "Cultivar IV Región"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 200
$ws.Range("N24").Value = 22000
$ws.Range("O24").Value = 22500
$ws.Range("P24").Value = 22250
$ws.Range("Q24").Value = "`$/bandeja 8 kilos"
$ws.Range("R24").Value = "Provincia de Limarí"
$ws.Range("S24").Value = 2781
$ws.Range("T24").Value = 8
